$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for two new "Hydro" columns between the existing Wind (F:G)
# and Plant (old H:I) / Infrastructure (old J:K) columns, without disturbing
# the column-width metadata tied to column positions 1-10.
#
# Old layout: A Country | B Elec | C Heat | D Solar% | E Solaryrs |
#             F Wind% | G Windyrs | H Plant% | I Plantyrs | J Infra% | K Infrayrs
# New layout: ... F Wind% | G Windyrs | H Hydro% | I Hydroyrs |
#             J Plant% | K Plantyrs | L Infra% | M Infrayrs

# Move the old "Plant"+"Infrastructure" block (H:K) two columns to the right (J:M)
$ws.Range("H1:K2").Copy()
$ws.Range("J1").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# The source block's H2:I2 values were copied forward to J2:K2 above; clear the
# stale leftovers in H2:I2 so the new Hydro data row can start blank.
$ws.Range("H2:I2").ClearContents()

# New "Hydro" header labels
$ws.Range("H1").Value = "Hydro interest rate"
$ws.Range("I1").Value = "Hydro lifetime (years)"

# Column I (9th col) widens slightly to fit the new "Hydro lifetime (years)" header
$ws.Columns.Item(9).ColumnWidth = 18.33

# Update the view: scroll right a bit and move the active selection
$null = $ws.Range("I10").Select()
$excel.ActiveWindow.ScrollColumn = 7
